$p = $ppt.ActivePresentation

# Slide 9 ("Installation Steps") - TextBox 4 (Shapes.Item(2)) holds the bullet
# list. The 3rd bullet ("Create S3 Bucket and Upload Deployment Artifacts")
# has its second run split into three runs so the text reads:
# "Create S3 " + "Bucket and Upload " + "AWS CloudFormation Deployment " + "Artifacts"
$s = $p.Slides.Item(9)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(3)
$run = $para.Runs(2)

$run.Text = "Bucket and Upload "
$run2 = $run.InsertAfter("AWS CloudFormation Deployment ")
$run2.InsertAfter("Artifacts")
